# Auto-generated edit script: add N7 policy rows (44-85) and related shared strings
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(44, 42, "I2_N7_T30_C140_0", 309.431978056256, 2.366000175476074, 0.9693721510508347, 33.86197805599477, 0, 0, 225.2600000002612, 50.31),
  @(45, 43, "I2_N7_T30_C140_0_P1", 1019.699639895617, 0.01200008392333984, 0, 26.35557902654179, 7.118526273374133, 5.792587142453418, 941.5899999999959, 53.08),
  @(46, 44, "I2_N7_T30_C140_0_P2", 396.5896102867711, 0.01600003242492676, 0, 27.5296102867676, 1.145862120850498, 1.145862120850498, 275.5500000000035, 93.50999999999999),
  @(47, 45, "I2_N7_T30_C140_0_P3", 396.5896102867718, 0.02200007438659668, 0, 27.5296102867676, 1.145862120850498, 1.145862120850498, 275.5500000000042, 93.50999999999999),
  @(48, 46, "I2_N7_T30_C140_0_P4", 357.599157859413, 0.01700019836425781, 0, 33.10915785941906, 3.9641518638922, 3.9641518638922, 264.119999999994, 60.37),
  @(49, 47, "I2_N7_T30_C140_0_P5", 1021.419132648, 0.009999990463256836, 0, 26.74913264799307, 4.450188684067006, 4.450188684067006, 393.770000000007, 600.9000000000001),
  @(50, 48, "I2_N7_T30_C140_0_P6", 696.9392770648811, 0.00800013542175293, 0, 36.70927706488092, 4.738617941974221, 4.738617941974221, 428.3500000000003, 231.88),
  @(51, 49, "I2_N7_T30_C210_0", 26.53458309510687, 0.9140000343322754, 0, 26.53458309510687, 0, 0, 0, 0),
  @(52, 50, "I2_N7_T30_C210_0_P1", 1091.229339571343, 0.009000062942504883, 0, 23.91796089943222, 4.77723927008741, 4.738617941974221, 620.2700000000234, 447.08),
  @(53, 51, "I2_N7_T30_C210_0_P2", 497.0523554984831, 0.01600003242492676, 0, 24.99064597944741, 6.065177049722176, 2.836886568752553, 366.5900000000053, 108.7),
  @(54, 52, "I2_N7_T30_C210_0_P3", 500.280645979443, 0.0130000114440918, 0, 24.99064597944741, 2.378003364169193, 2.378003364169193, 422.4199999999955, 52.87),
  @(55, 53, "I2_N7_T30_C210_0_P4", 135.0555790265418, 0.006000041961669922, 0, 26.35557902654179, 0, 0, 0, 108.7),
  @(56, 54, "I2_N7_T30_C210_0_P5", 1091.43471156119, 0.006999969482421875, 0, 24.3115145208835, 3.366602695185008, 3.139799735480239, 969.9100000000115, 97.44),
  @(57, 55, "I2_N7_T30_C210_0_P6", 147.3892770648809, 0, 0, 36.70927706488092, 0, 0, 0, 110.68),
  @(58, 56, "I2_N7_T30_C280_0", 26.29875013250022, 1.130999803543091, 0, 26.29875013250022, 0, 0, 0, 0),
  @(59, 57, "I2_N7_T30_C280_0_P1", 1676.610268005549, 0.006999969482421875, [double]"7.11293361003758E-06", 21.92482366701793, 4.943173603449611, 4.738617941974221, 1347.510000000006, 307.3800000000001),
  @(60, 58, "I2_N7_T30_C280_0_P2", 1014.423453125771, 0.006000041961669922, 0, 22.00055025312762, 5.894019291323108, 3.646922163947599, 793.1700000000191, 201.5),
  @(61, 59, "I2_N7_T30_C280_0_P3", 391.2679608994313, 0.009999990463256836, 0, 23.91796089943222, 3.9641518638922, 3.9641518638922, 218.7199999999991, 148.63),
  @(62, 60, "I2_N7_T30_C280_0_P4", 77.8606459794474, 0.0149998664855957, 0, 24.99064597944741, 0, 0, 0, 52.87),
  @(63, 61, "I2_N7_T30_C280_0_P5", 1674.399505879263, 0.003999948501586914, 0, 23.33006819063954, 6.41079396126198, 2.590231649872266, 1129.060000000013, 525.83),
  @(64, 62, "I2_N7_T30_C280_0_P6", 147.3892770648809, 0.004999876022338867, 0, 36.70927706488092, 0, 0, 0, 110.68),
  @(65, 63, "I2_N7_T100_C140_0", 307.6174397573979, 181.175999879837, 0, 32.04743975738332, 0, 0, 275.5700000000146, 0),
  @(66, 64, "I2_N7_T100_C140_0_P1", 1017.797288545571, 0.01200008392333984, 0, 26.35557902654179, 8.443180413891367, 5.214889932921746, 691.7199999999982, 302.95),
  @(67, 65, "I2_N7_T100_C140_0_P2", 396.5896102867676, 0.01099991798400879, 0, 27.5296102867676, 1.145862120850498, 1.145862120850498, 275.55, 93.50999999999999),
  @(68, 66, "I2_N7_T100_C140_0_P3", 396.5896102867676, 0.01200008392333984, 0, 27.5296102867676, 1.145862120850498, 1.145862120850498, 275.55, 93.50999999999999),
  @(69, 67, "I2_N7_T100_C140_0_P4", 357.4462074868381, 0.01200008392333984, 0, 32.95620748683876, 3.9641518638922, 3.9641518638922, 271.4099999999993, 53.08),
  @(70, 68, "I2_N7_T100_C140_0_P5", 1021.419132647999, 0.009999990463256836, 0, 26.74913264799307, 4.450188684067006, 4.450188684067006, 645.0300000000061, 349.64),
  @(71, 69, "I2_N7_T100_C140_0_P6", 651.7408010327504, 0.01200008392333984, 0, 34.03080103274972, 2.843170765184533, 2.843170765184533, 448.6400000000006, 169.07),
  @(72, 70, "I2_N7_T100_C210_0", 26.35557902654179, 5.176999807357788, 0, 26.3555790265418, 0, 0, 0, 0),
  @(73, 71, "I2_N7_T100_C210_0_P1", 1091.229339571329, 0.00800013542175293, 0, 23.91796089943222, 4.77723927008741, 4.738617941974221, 717.7100000000102, 349.64),
  @(74, 72, "I2_N7_T100_C210_0_P2", 500.2806459794576, 0.006999969482421875, 0, 24.99064597944741, 2.378003364169193, 2.378003364169193, 235.9700000000101, 239.32),
  @(75, 73, "I2_N7_T100_C210_0_P3", 500.2806459794576, 0.008999824523925781, 0, 24.99064597944741, 2.378003364169193, 2.378003364169193, 235.9700000000101, 239.32),
  @(76, 74, "I2_N7_T100_C210_0_P4", 26.35557902654179, 0.004999876022338867, 0, 26.35557902654179, 0, 0, 0, 0),
  @(77, 75, "I2_N7_T100_C210_0_P5", 1091.321310081332, 0.009000062942504883, 0, 24.3115145208835, 3.754788217841379, 3.414583778284226, 1014.270000000005, 53.08),
  @(78, 76, "I2_N7_T100_C210_0_P6", 142.7308010327497, 0.004999876022338867, 0, 34.03080103274972, 0, 0, 0, 108.7),
  @(79, 77, "I2_N7_T100_C280_0", 24.99064597944741, 5.861999988555908, 0, 24.99064597944741, 0, 0, 0, 0),
  @(80, 78, "I2_N7_T100_C280_0_P1", 1676.610268005544, 0.007999897003173828, [double]"6.830343646905243E-06", 21.92482366701793, 4.943173603449611, 4.738617941974221, 1436.870000000001, 218.02),
  @(81, 79, "I2_N7_T100_C280_0_P2", 1013.861678843903, 0.007999897003173828, 0, 22.00055025312762, 6.773023273111587, 3.9641518638922, 846.039999999995, 148.63),
  @(82, 80, "I2_N7_T100_C280_0_P3", 391.2525123681869, 0.009000062942504883, 0, 23.91796089943222, 2.598412980545262, 2.582964449299987, 367.35, 0),
  @(83, 81, "I2_N7_T100_C280_0_P4", 24.99064597944741, 0.006000041961669922, 0, 24.99064597944741, 0, 0, 0, 0),
  @(84, 82, "I2_N7_T100_C280_0_P5", 1674.399505879249, 0.007999897003173828, 0, 23.33006819063954, 6.41079396126198, 2.590231649872266, 1384.939999999999, 269.95),
  @(85, 83, "I2_N7_T100_C280_0_P6", 142.7308010327497, 0.006000041961669922, 0, 34.03080103274972, 0, 0, 0, 108.7)
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value = $r[1]
  $ws.Cells.Item($rowNum, 2).Value = $r[2]
  for ($col = 3; $col -le 10; $col++) {
    $ws.Cells.Item($rowNum, $col).Value = $r[$col]
  }
}

# Copy the bold/bordered style used in column A (row 43) onto the new column-A cells
$ws.Range("A43").Copy()
$ws.Range("A44:A85").PasteSpecial(-4122)
$excel.CutCopyMode = 0
